# 03.02.2025 - Poprawa błędu podczas importowania raportowania zleceń
# Adds two new "to do" rows at the bottom of the Sheet1 task list and
# scrolls/selects the sheet roughly where the edit was made.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- New row 62 ---------------------------------------------------------
$ws.Cells.Item(62, 1).Value = 61
$ws.Cells.Item(62, 2).Value = "Wyliczenia"
$ws.Cells.Item(62, 3).Value = "Dodać zabezpieczenie przed powtórnym zapisem. Przycisk jest już wstawiony"
$ws.Cells.Item(62, 4).Value = 0

# --- New row 63 ---------------------------------------------------------
$ws.Cells.Item(63, 1).Value = 62
$ws.Cells.Item(63, 2).Value = "Ustawienia"
$ws.Cells.Item(63, 3).Value = "Dodać listę użytkowników w bazie danych.W tej chwili jest wpisana na sztywno w kodzie. "
$ws.Cells.Item(63, 4).Value = 0

# Match the formatting used by the rest of the table:
#   column A/B -> vertical-top, column C -> wrap text, column D -> vertical-center
$rowsToFormat = @(62, 63)
foreach ($r in $rowsToFormat) {
    $ws.Cells.Item($r, 1).VerticalAlignment = -4160   # xlVAlignTop
    $ws.Cells.Item($r, 2).VerticalAlignment = -4160   # xlVAlignTop
    $ws.Cells.Item($r, 3).WrapText = $true
    $ws.Cells.Item($r, 4).VerticalAlignment = -4108   # xlVAlignCenter
}

# Leave the sheet scrolled/selected near the edited rows, matching where the
# author was working when the file was saved.
$ws.Range("C60:D61").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 46
try { $excel.ActiveWindow.Top = -16320 } catch {}
